# Refresh the cryptos listing: update Price (col D) and Volume(1h) (col E)
# for each coin row, and re-sync rows 29-31 whose coins were re-ranked.
#
# Several "Price" values are strings that happen to look numeric (e.g.
# "0.999", "26.33"). To keep them stored as text (matching the original
# inlineStr cells) instead of being auto-converted to numbers, they are
# written with a leading apostrophe (forces text entry) and the cell
# style is immediately reset to "Normal" so no stray number-format style
# is left attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.780.58'
$ws.Range("E2").Value = '  -2.40%  '
$ws.Range("D3").Value = '2.495.89'
$ws.Range("E3").Value = '  -4.50%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = "'554.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.06%  '
$ws.Range("D6").Value = "'146.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.26%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").Value = "'0.604"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.74%  '
$ws.Range("D9").Value = '2.491.60'
$ws.Range("E9").Value = '  -4.57%  '
$ws.Range("E10").Value = '  -6.72%  '
$ws.Range("D11").Value = "'5.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.09%  '
$ws.Range("D13").Value = "'0.360"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.72%  '
$ws.Range("D14").Value = "'26.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.15%  '
$ws.Range("D15").Value = '2.939.43'
$ws.Range("E15").Value = '  -4.71%  '
$ws.Range("D16").Value = "'0.0000168"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -7.38%  '
$ws.Range("D17").Value = '61.678.52'
$ws.Range("E17").Value = '  -2.51%  '
$ws.Range("D18").Value = '2.498.23'
$ws.Range("E18").Value = '  -4.26%  '
$ws.Range("D19").Value = "'11.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.50%  '
$ws.Range("D20").Value = "'7.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.85%  '
$ws.Range("E21").Value = '  -5.77%  '
$ws.Range("D22").Value = "'323.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.25%  '
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("E24").Value = '  -1.99%  '
$ws.Range("D25").Value = "'64.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.26%  '
$ws.Range("D26").Value = '0.0₃0999'
$ws.Range("E26").Value = '  -5.36%  '
$ws.Range("D27").Value = '2.601.85'
$ws.Range("E27").Value = '  -5.12%  '
$ws.Range("D28").Value = "'1.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.56%  '
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").Value = "'8.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -6.73%  '
$ws.Range("B30").Value = 'Bittensor'
$ws.Range("C30").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D30").Value = "'539.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.47%  '
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").Value = "'0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.18%  '
$ws.Range("D32").Value = "'7.63"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.58%  '
$ws.Range("E33").Value = '  -5.17%  '
$ws.Range("D34").Value = "'1.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.14%  '
$ws.Range("D35").Value = "'1.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.38%  '
$ws.Range("D36").Value = "'6.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.28%  '
$ws.Range("D37").Value = "'4.94"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.76%  '
$ws.Range("D38").Value = "'0.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("D39").Value = "'0.386"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'18.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.08%  '
$ws.Range("D41").Value = "'148.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.21%  '
$ws.Range("D42").Value = "'1.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.03%  '
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").Value = "'40.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.65%  '
$ws.Range("D45").Value = "'2.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.19%  '
$ws.Range("D46").Value = "'148.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.72%  '
$ws.Range("D47").Value = "'3.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.87%  '
$ws.Range("D48").Value = "'20.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -12.11%  '
$ws.Range("D49").Value = "'0.0537"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -7.72%  '
$ws.Range("D50").Value = "'0.601"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.85%  '
$ws.Range("D51").Value = "'0.0950"
$ws.Range("D51").Style = "Normal"
